# Updated cryptos list on Thu Mar 16 14:41:57 UTC 2023 with GitHub Actions
#
# Rewrites the Price (D) and Volume(1h) (E) columns for the crypto table
# with freshly scraped figures. A couple of coins (EnergySwap /
# PancakeSwap) also swapped table position, so those two rows get their
# Coin name + Link + Price + Volume replaced wholesale.
#
# Price/Volume values are plain text in the sheet (e.g. "1.190", "0.9993",
# "  -0.60%  ") even though many of them look like numbers. Assigning such
# a string straight to Range.Value makes Excel silently reinterpret it as
# a real number (dropping the trailing zero / exact text), and flips the
# cell to quote-prefixed text if we try to force it with a leading
# apostrophe - neither matches the source data. Temporarily forcing the
# whole target range to Text format before writing keeps every value as
# literal text; re-applying the Normal style afterwards restores the
# original (unstyled) cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.888.30"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.658.25"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "323.39"
$ws.Range("E5").Value = "  +4.54%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "0.3636"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").Value = "47.52"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "0.3262"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "1.133"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").Value = "0.07071"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "6.054"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "19.54"
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").Value = "1.657.80"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "0.00001048"
$ws.Range("E17").Value = "  -3.66%  "
$ws.Range("D18").Value = "0.06591"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "0.9993"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "79.00"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("D21").Value = "5.911"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").Value = "15.75"
$ws.Range("E22").Value = "  -5.30%  "
$ws.Range("D23").Value = "12.68"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").Value = "24.858.81"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "2.436"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "2.456"
$ws.Range("E26").Value = "  -7.07%  "
$ws.Range("D27").Value = "147.63"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "18.62"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("D29").Value = "1.839.14"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").Value = "1.199"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "125.12"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").Value = "4.097"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").Value = "5.757"
$ws.Range("E33").Value = "  -9.38%  "
$ws.Range("D34").Value = "0.08443"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "1.642"
$ws.Range("E35").Value = "  -5.18%  "
$ws.Range("D36").Value = "12.25"
$ws.Range("E36").Value = "  -7.47%  "
$ws.Range("D37").Value = "1.280"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "5.160"
$ws.Range("E38").Value = "  -3.43%  "
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").Value = "0.06060"
$ws.Range("E40").Value = "  -5.07%  "
$ws.Range("D41").Value = "8.353"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("D42").Value = "0.2068"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "0.5935"
$ws.Range("E44").Value = "  -4.04%  "

# EnergySwap and PancakeSwap swapped ranking positions this run.
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.51"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.885"
$ws.Range("E46").Value = "  +2.66%  "

$ws.Range("D47").Value = "0.5613"
$ws.Range("E47").Value = "  -4.66%  "
$ws.Range("D48").Value = "124.89"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "1.946"
$ws.Range("E49").Value = "  -3.57%  "
$ws.Range("D50").Value = "0.06980"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("E51").Value = "  -0.04%  "

# Restore the original (unstyled) formatting now that every value is
# committed as literal text.
$priceVolumeRange.Style = "Normal"
